# rdbms-02.data.xlsx -- add a new "customers SQL file" property (row 10, which
# was previously blank) and append a new "nexial.verbose" property row (row 11)
# to the DynamicSQL3 sheet.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("DynamicSQL3")
$ws0 = $wb.Worksheets.Item(1)

$ws.Activate()

# ---- row 10: fill in the previously-empty property row ----------------
$ws.Range("A10").Value = "customers SQL file"
$ws.Range("B10").Value = '$(syspath|data|fullpath)/rdbms-02-queries-to-csv.sql'

# ---- row 11: new property, styled like the bordered "nexial.verbose"
#      row that already exists at the bottom of the #default sheet -----
$ws0.Range("A1:B1").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A11").Value = "nexial.verbose"

# write "false" as literal text (not the Boolean FALSE) by routing it
# through a formula-to-value paste, reusing the existing shared string
$ws.Range("B11").Formula = '="false"'
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)       # xlPasteValues

$ws.Range("A11:B11").Select()
